$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.472.24"
$ws.Range("E2").Value = "  -1.40%  "

$ws.Range("D3").Value = "1.747.00"
$ws.Range("E3").Value = "  -1.23%  "

$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").Value = "'328.15"
$ws.Range("E5").Value = "  +0.28%  "

$ws.Range("D6").Value = "'1.004"
$ws.Range("E6").Value = "  -0.04%  "

$ws.Range("D7").Value = "'0.4717"
$ws.Range("E7").Value = "  +5.18%  "

$ws.Range("D8").Value = "'0.3519"
$ws.Range("E8").Value = "  -1.62%  "

$ws.Range("D9").Value = "'42.74"
$ws.Range("E9").Value = "  +0.98%  "

$ws.Range("D10").Value = "'0.07433"
$ws.Range("E10").Value = "  -0.34%  "

$ws.Range("D11").Value = "'1.072"
$ws.Range("E11").Value = "  -2.30%  "

$ws.Range("D12").Value = "'1.004"
$ws.Range("E12").Value = "  -0.15%  "

$ws.Range("D13").Value = "'20.28"
$ws.Range("E13").Value = "  -2.89%  "

$ws.Range("D14").Value = "'6.011"
$ws.Range("E14").Value = "  -0.67%  "

$ws.Range("D15").Value = "1.745.80"
$ws.Range("E15").Value = "  -1.31%  "

$ws.Range("D16").Value = "'7.025"
$ws.Range("E16").Value = "  -2.60%  "

$ws.Range("D17").Value = "'91.69"
$ws.Range("E17").Value = "  -1.59%  "

$ws.Range("E18").Value = "  +0.93%  "

$ws.Range("D19").Value = "'0.06368"
$ws.Range("E19").Value = "  -1.30%  "

$ws.Range("E20").Value = "  -0.03%  "

$ws.Range("E21").Value = "  -3.26%  "

$ws.Range("D22").Value = "'5.745"
$ws.Range("E22").Value = "  -1.23%  "

$ws.Range("D23").Value = "27.534.47"
$ws.Range("E23").Value = "  -1.37%  "

$ws.Range("E24").Value = "  -3.19%  "

$ws.Range("D25").Value = "'2.148"
$ws.Range("E25").Value = "  +1.77%  "

$ws.Range("D26").Value = "'161.00"
$ws.Range("E26").Value = "  -1.39%  "

$ws.Range("D27").Value = "'19.82"
$ws.Range("E27").Value = "  -2.38%  "

$ws.Range("D28").Value = "1.947.04"
$ws.Range("E28").Value = "  -1.90%  "

$ws.Range("D29").Value = "'2.181"
$ws.Range("E29").Value = "  -1.22%  "

$ws.Range("D30").Value = "'121.39"
$ws.Range("E30").Value = "  -3.42%  "

$ws.Range("D31").Value = "'1.044"
$ws.Range("E31").Value = "  -5.28%  "

$ws.Range("D32").Value = "'0.09300"
$ws.Range("E32").Value = "  +1.96%  "

$ws.Range("E33").Value = "  -0.57%  "

$ws.Range("D34").Value = "'5.456"
$ws.Range("E34").Value = "  -1.99%  "

$ws.Range("D35").Value = "'0.02240"
$ws.Range("E35").Value = "  -2.30%  "

$ws.Range("D36").Value = "'11.32"
$ws.Range("E36").Value = "  -4.93%  "

$ws.Range("D37").Value = "'0.05927"
$ws.Range("E37").Value = "  -3.02%  "

$ws.Range("E38").Value = "  -2.33%  "

$ws.Range("D39").Value = "'4.832"
$ws.Range("E39").Value = "  -2.64%  "

$ws.Range("E40").Value = "  +2.39%  "

$ws.Range("D41").Value = "'0.6066"
$ws.Range("E41").Value = "  -4.43%  "

$ws.Range("D42").Value = "'1.170"
$ws.Range("E42").Value = "  -1.27%  "

$ws.Range("D43").Value = "'7.703"
$ws.Range("E43").Value = "  -2.84%  "

$ws.Range("D44").Value = "'3.725"
$ws.Range("E44").Value = "  -0.34%  "

$ws.Range("D45").Value = "'12.86"
$ws.Range("E45").Value = "  -3.50%  "

$ws.Range("D46").Value = "'0.5687"
$ws.Range("E46").Value = "  -3.72%  "

$ws.Range("D47").Value = "'122.62"
$ws.Range("E47").Value = "  -0.19%  "

$ws.Range("D48").Value = "'1.900"
$ws.Range("E48").Value = "  -3.00%  "

$ws.Range("D49").Value = "'1.129"
$ws.Range("E49").Value = "  -0.79%  "

$ws.Range("D50").Value = "'0.06751"
$ws.Range("E50").Value = "  -2.50%  "

$ws.Range("E51").Value = "  -2.50%  "
